$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8511
    $ws.Range("F3").Value = 193
    $ws.Range("F4").Value = 379
    $ws.Range("F5").Value = 23
}
